$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 5-10, column C: give the already-present Non-hazardous/Hazardous
#    label cells an explicit (but visually neutral) font touch, matching the
#    "applyFont" style the author picked up after re-saving with Excel.
# ---------------------------------------------------------------------------
foreach ($addr in @("C5", "C6", "C7", "C8", "C9", "C10")) {
    $ws.Range($addr).Font.Name = "Calibri"
}

# ---------------------------------------------------------------------------
# 2) New EURAL rows for the "chapter 04" keyflow (leather/textile industry
#    waste), appended below the existing chapter-20 rows.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 18; Code = "040102"; Name = "loogafval"; Haz = "Non-hazardous" },
    @{ Row = 19; Code = "040106"; Name = "chroomhoudend slib, met name van afvalwaterbehandeling ter plaatse"; Haz = "Non-hazardous" },
    @{ Row = 20; Code = "040109"; Name = "afval van bewerking en afwerking"; Haz = "Non-hazardous" },
    @{ Row = 21; Code = "040209"; Name = "afval van composietmaterialen (geïmpregneerde textiel, elastomeren, plastomeren)"; Haz = "Non-hazardous" },
    @{ Row = 22; Code = "040214"; Name = "afval van afwerking dat organische oplosmiddelen bevat"; Haz = "Hazardous" },
    @{ Row = 23; Code = "040216"; Name = "kleurstoffen en pigmenten die Hazardouse stoffen bevatten"; Haz = "Hazardous" },
    @{ Row = 24; Code = "040219"; Name = "slib van afvalwaterbehandeling ter plaatse dat Hazardouse stoffen bevat"; Haz = "Hazardous" },
    @{ Row = 25; Code = "040220"; Name = "niet onder 04 02 19 vallend slib van afvalwaterbehandeling ter plaatse"; Haz = "Non-hazardous" },
    @{ Row = 26; Code = "040222"; Name = "afval van verwerkte textielvezels"; Haz = "Non-hazardous" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $codeCell = $ws.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.HorizontalAlignment = -4131
    $codeCell.Value = $item.Code

    $nameCell = $ws.Cells.Item($r, 2)
    $nameCell.IndentLevel = 1
    $nameCell.Value = $item.Name

    $hazCell = $ws.Cells.Item($r, 3)
    $hazCell.Value = $item.Haz
}

# ---------------------------------------------------------------------------
# 3) Selection cursor moves to C17 after the edit.
# ---------------------------------------------------------------------------
$ws.Range("C17").Select()
